$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select rows 2 and 3 (the empty START_SCENE/INTRO ... END_SCENE block)
# before removing them, mirroring how this edit was made interactively.
$ws.Range("A2:XFD3").Select()

# Delete rows 2 and 3, shifting all following rows up by two.
$ws.Rows("2:3").Delete()

# Update the scene counter in A1 (was counting 4 scenes, now 3).
$ws.Range("A1").Value = 3
